# Weekly data refresh: a new price record (2023-12-04) is added for
# "Arveja Verde" at row 84, pushing all subsequent records down by one row.
# The last existing record (old row 167) ends up duplicated as the new
# final row (168), since nothing is removed - everything simply shifts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 84; Excel shifts rows 84..167 down to 85..168
# and the sheet dimension grows from A1:R167 to A1:R168 automatically.
$ws.Rows("84").Insert()

# Populate the newly inserted row 84 with the new record's data.
$ws.Range("A84").Value = 5
$ws.Range("B84").Value = "Macroferia Regional de Talca"
$ws.Range("C84").Value = "Maule"
$ws.Range("D84").Value = 45264
$ws.Range("E84").Value = 7
$ws.Range("F84").Value = 100112022
$ws.Range("G84").Value = "Arveja Verde"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 400
$ws.Range("K84").Value = 23000
$ws.Range("L84").Value = 23000
$ws.Range("M84").Value = 23000
$ws.Range("N84").Value = "`$/saco 25 kilos"
$ws.Range("O84").Value = "Región del Maule"
$ws.Range("P84").Value = 920
$ws.Range("Q84").Value = 25
$ws.Range("R84").Value = "Hortaliza"
